$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- Phase 0: reset pre-existing styles on D3:D12 so formula-setting does not inherit them ----
$ws.Range("D3:D12").Style = "Normal"

# ---- Phase 1: set plain values for columns A,B,C,D,F,G ----
$ws.Range("A2").Value = [double]"3.16"
$ws.Range("B2").Value = [double]"5"
$ws.Range("C2").Value = [double]"3.16"
$ws.Range("D2").Value = [double]"0.19400000000000001"
$ws.Range("A3").Value = [double]"3.16"
$ws.Range("B3").Value = [double]"10"
$ws.Range("C3").Value = [double]"3.16"
$ws.Range("D3").Value = [double]"9.8400000000000001E-2"
$ws.Range("A4").Value = [double]"3.16"
$ws.Range("B4").Value = [double]"20"
$ws.Range("C4").Value = [double]"3.16"
$ws.Range("D4").Value = [double]"0.05"
$ws.Range("A5").Value = [double]"3.16"
$ws.Range("B5").Value = [double]"40"
$ws.Range("C5").Value = [double]"3.16"
$ws.Range("D5").Value = [double]"0"
$ws.Range("A6").Value = [double]"3.16"
$ws.Range("B6").Value = [double]"80.599999999999994"
$ws.Range("C6").Value = [double]"3.16"
$ws.Range("D6").Value = [double]"0"
$ws.Range("A7").Value = [double]"3.16"
$ws.Range("B7").Value = [double]"800"
$ws.Range("C7").Value = [double]"3.12"
$ws.Range("D7").Value = [double]"1.2E-5"
$ws.Range("A8").Value = [double]"3.16"
$ws.Range("B8").Value = [double]"1600"
$ws.Range("C8").Value = [double]"3.12"
$ws.Range("D8").Value = [double]"1.5999999999999999E-5"
$ws.Range("A9").Value = [double]"3.16"
$ws.Range("B9").Value = [double]"3200"
$ws.Range("C9").Value = [double]"3.06"
$ws.Range("D9").Value = [double]"1.7E-5"
$ws.Range("A10").Value = [double]"3.16"
$ws.Range("B10").Value = [double]"6480"
$ws.Range("C10").Value = [double]"2.84"
$ws.Range("D10").Value = [double]"1.56E-5"
$ws.Range("A11").Value = [double]"3.16"
$ws.Range("B11").Value = [double]"12780"
$ws.Range("C11").Value = [double]"2.2000000000000002"
$ws.Range("D11").Value = [double]"1.4E-5"
$ws.Range("A12").Value = [double]"3.16"
$ws.Range("B12").Value = [double]"18070"
$ws.Range("C12").Value = [double]"1.68"
$ws.Range("D12").Value = [double]"1.29E-5"
$ws.Range("A13").Value = [double]"3.16"
$ws.Range("B13").Value = [double]"22500"
$ws.Range("C13").Value = [double]"1.3"
$ws.Range("D13").Value = [double]"1.1600000000000001E-5"
$ws.Range("A14").Value = [double]"3.16"
$ws.Range("B14").Value = [double]"30000"
$ws.Range("C14").Value = [double]"0.84"
$ws.Range("D14").Value = [double]"1.08E-5"
$ws.Range("A15").Value = [double]"2.82"
$ws.Range("B15").Value = [double]"38100"
$ws.Range("C15").Value = [double]"0.52800000000000002"
$ws.Range("D15").Value = [double]"9.0000000000000002E-6"
$ws.Range("A16").Value = [double]"2.78"
$ws.Range("B16").Value = [double]"50000"
$ws.Range("C16").Value = [double]"0.3"
$ws.Range("D16").Value = [double]"8.2800000000000003E-6"
$ws.Range("A17").Value = [double]"0.27600000000000002"
$ws.Range("B17").Value = [double]"57070"
$ws.Range("C17").Value = [double]"2.52E-2"
$ws.Range("F17").Value = [double]"0.28749999999999998"
$ws.Range("G17").Value = [double]"0.17499999999999999"
$ws.Range("A18").Value = [double]"2.7399999999999998E-3"
$ws.Range("B18").Value = [double]"80000"
$ws.Range("C18").Value = [double]"1.26E-2"
$ws.Range("F18").Value = [double]"0.2"
$ws.Range("G18").Value = [double]"3.2000000000000001E-2"
$ws.Range("A19").Value = [double]"0.16500000000000001"
$ws.Range("B19").Value = [double]"191100"
$ws.Range("C19").Value = [double]"1.7600000000000001E-2"
$ws.Range("F19").Value = [double]"1"
$ws.Range("G19").Value = [double]"0"
$ws.Range("A20").Value = [double]"0.16700000000000001"
$ws.Range("B20").Value = [double]"1000000"
$ws.Range("C20").Value = [double]"6.4799999999999996E-2"
$ws.Range("F20").Value = [double]"8.3000000000000007"
$ws.Range("G20").Value = [double]"5"
$ws.Range("A21").Value = [double]"0.17"
$ws.Range("B21").Value = [double]"2700000"
$ws.Range("C21").Value = [double]"0.08"
$ws.Range("F21").Value = [double]"8.5"
$ws.Range("G21").Value = [double]"1.6"

# ---- Phase 2: set formulas ----
# Group E11:E16 first (-> si=1), then E3:E10 (-> si=0), matching target shared indices
$ws.Range("E11:E16").Formula = "=D11*B11*360"
$ws.Range("E3:E10").Formula = "=D3*B3*360"
$ws.Range("E2").Formula = "=D2*B2*360"
$ws.Range("E17").Formula = "=180-ASIN(G17/F17)"
$ws.Range("E18").Formula = "=180-ASIN(G18/F18)"
$ws.Range("E19").Formula = "=180-ASIN(G19/F19)"
$ws.Range("E20").Formula = "=ASIN(G20/F20)"
$ws.Range("E21").Formula = "=ASIN(G21/F21)"

# ---- Phase 3: apply final styles (numFmt scientific 0.00E+00 = existing style index 1) ----
$ws.Range("A17:A21").NumberFormat = "0.00E+00"
$ws.Range("B10:B21").NumberFormat = "0.00E+00"
$ws.Range("C14:C21").NumberFormat = "0.00E+00"
$ws.Range("D3:D17").NumberFormat = "0.00E+00"
$ws.Range("F17:F21").NumberFormat = "0.00E+00"
$ws.Range("E19").NumberFormat = "0.00E+00"

# ---- Phase 4: header cells F1/G1 (new shared strings "A" and "C") ----
$ws.Range("F1").Value = "A"
$ws.Range("G1").Value = "C"

# ---- Phase 5: sheet/view selections ----
$ws.Activate()
$ws.Range("B27").Select() | Out-Null
